$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '250.56'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '24.27'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.955'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05907'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.425'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.336'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.7959'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1489'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07783'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03307'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09244'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.555'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001669'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04754'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.006223'
$ws.Range('E18').Value = '17TigerCashTCH'
$ws.Range('B19').Value = 'HotbitToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.005582'
$ws.Range('E19').Value = '18HotbitTokenHTB'
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.001065'
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('B21').Value = 'NitroEx'
$ws.Range('C21').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0001498'
$ws.Range('E21').Value = '20NitroExNTX'
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.710'
$ws.Range('E22').Value = '21LEOLEO'
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.211'
$ws.Range('E23').Value = '22BTSETokenBTSE'
$ws.Range('B24').Value = 'One'
$ws.Range('C24').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.01271'
$ws.Range('E24').Value = '23OneONEBestin24h'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.3353'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1252'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0006468'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04400'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007029'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01003'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.002458'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00005891'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9891'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1103'
$ws.Range('E49').Value = '48BOLOBOLOWorstin24h'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002098'
